$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "logCompare" between "stress" and "status".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "logCompare"
$newSheet.Move($wb.Worksheets.Item("status"))

# Re-fetch by name - the handle obtained before Move() goes stale afterwards.
$logCompare = $wb.Worksheets.Item("logCompare")

# ---------------------------------------------------------------------------
# 2. Populate "logCompare" - header row uses the same text-format style as
#    the other sheets' header rows.
# ---------------------------------------------------------------------------
$logCompare.Range("A1:C1").NumberFormat = "@"
$logCompare.Range("A1").Value = "Scenario No."
$logCompare.Range("B1").Value = "key"
$logCompare.Range("C1").Value = "value"

# Data rows 2-11 (scenario1), filled row by row, column by column so the
# shared-string table is populated in the same order the source data was.
$logCompare.Range("A2").Value = "scenario1"
$logCompare.Range("B2").Value = "match"
$logCompare.Range("C2").Value = "F41F42"

$logCompare.Range("A3").Value = "scenario1"
$logCompare.Range("B3").Value = "comp"
$logCompare.Range("C3").Value = "011GT000006"

$logCompare.Range("A4").Value = "scenario1"
$logCompare.Range("B4").Value = "exclude"
$logCompare.Range("C4").Value = 7011
$logCompare.Range("C4").NumberFormat = "#,##0"

$logCompare.Range("A5").Value = "scenario1"
$logCompare.Range("B5").Value = "uatLog"
$logCompare.Range("C5").Value = "C:\\FINsim\\UAT.mlg"

$logCompare.Range("A6").Value = "scenario1"
$logCompare.Range("B6").Value = "productionLog"
$logCompare.Range("C6").Value = "C:\\FINsim\\production.mlg"

$logCompare.Range("A7").Value = "scenario1"
$logCompare.Range("B7").Value = "deviceName"
$logCompare.Range("C7").Value = "device_name_1"

$logCompare.Range("A8").Value = "scenario1"
$logCompare.Range("B8").Value = "reportName"
$logCompare.Range("C8").Value = "reportName_1.rtf"

$logCompare.Range("A9").Value = "scenario1"
$logCompare.Range("B9").Value = "location"
$logCompare.Range("C9").Value = "C:\\FINsim\\"

$logCompare.Range("A10").Value = "scenario1"
$logCompare.Range("B10").Value = "sprAddress"
$logCompare.Range("C10").Value = "127.0.0.3"

$logCompare.Range("A11").Value = "scenario1"
$logCompare.Range("B11").Value = "sprPort"
$logCompare.Range("C11").Value = 1721

# Row 13 (scenario2) is filled before row 12 - matches the shared-string
# ordering recorded in the target workbook.
$logCompare.Range("A13").Value = "scenario2"
$logCompare.Range("B13").Value = "test"
$logCompare.Range("C13").Value = "testValue"

$logCompare.Range("A12").Value = "scenario2"
$logCompare.Range("B12").Value = "match"
$logCompare.Range("C12").Value = "F43"

# ---------------------------------------------------------------------------
# 3. Cosmetics on "logCompare": column widths + final selection/active cell.
# ---------------------------------------------------------------------------
$logCompare.Columns.Item(1).ColumnWidth = 17.43
$logCompare.Columns.Item(2).ColumnWidth = 18.57
$logCompare.Columns.Item(3).ColumnWidth = 33.14

# ---------------------------------------------------------------------------
# 4. Update the "stress" sheet selection (no longer the active tab).
# ---------------------------------------------------------------------------
$stress = $wb.Worksheets.Item("stress")
$stress.Range("A1:C1").Select()

# ---------------------------------------------------------------------------
# 5. Make "logCompare" the active sheet/tab, with G13 selected - this also
#    drives workbookView's activeTab and the per-sheet tabSelected flag.
# ---------------------------------------------------------------------------
$logCompare.Select()
$logCompare.Range("G13").Select()
